{"js": "// Apply strikethrough formatting to the \"Kubernetes: ...\" list item\n// (both the run text and the paragraph mark), matching the target diff.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.indexOf(\"Kubernetes:\") !== -1);\nif (!target) {\n  throw new Error('Paragraph containing \"Kubernetes:\" not found');\n}\n\n// Setting strikeThrough on the paragraph's font applies it to the whole\n// paragraph range, including the trailing paragraph mark (pPr/rPr), which\n// mirrors selecting the full line (pilcrow included) and toggling\n// Strikethrough in the Word UI.\ntarget.font.strikeThrough = true;\n\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to the \"Kubernetes: ...\" list item\n# (both the run text and the paragraph mark), matching the target diff.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"Kubernetes:\")\n\nif (-not $rng.Find.Found) {\n    throw \"Could not find paragraph containing 'Kubernetes:'\"\n}\n\n# Grab the whole paragraph (Paragraph.Range includes the trailing paragraph\n# mark), so toggling StrikeThrough affects both the run and the pPr/rPr,\n# mirroring selecting the full line including the pilcrow in the Word UI.\n$para = $rng.Paragraphs(1)\n$prange = $para.Range\n$prange.Font.StrikeThrough = 1\n"}
